# Auto-generated Excel COM-interop script to apply the diff to before.xlsx
# Updates currentAveragePrice / Leve price / profit columns (H-N) for specific rows
# across the ALC, ARM, BSM, CUL, GSM, LTW, WVR sheets, matching the scheduled-runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 571.0714
$ws.Range("I28").Value = 571.0714
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 571.0714
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -86.07140000000004
$ws.Range("N28").ClearContents()

$ws.Range("H33").Value = 100.69231
$ws.Range("I33").Value = 67.416664
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 67.416664
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 161.583336
$ws.Range("N33").Value = -958

$ws.Range("H40").Value = 3381.0625
$ws.Range("I40").Value = 5419.6
$ws.Range("J40").Value = 2454.4546
$ws.Range("K40").Value = 5419.6
$ws.Range("L40").Value = 2454.4546
$ws.Range("M40").Value = -5244.6
$ws.Range("N40").Value = -2804.4546

$ws.Range("H98").Value = 1348.0526
$ws.Range("I98").Value = 1461.5714
$ws.Range("J98").Value = 1030.2
$ws.Range("K98").Value = 1461.5714
$ws.Range("L98").Value = 1030.2
$ws.Range("M98").Value = 36.42859999999996
$ws.Range("N98").Value = -4026.2

$ws.Range("H122").Value = 1348.0526
$ws.Range("I122").Value = 1461.5714
$ws.Range("J122").Value = 1030.2
$ws.Range("K122").Value = 4384.7142
$ws.Range("L122").Value = 3090.6
$ws.Range("M122").Value = -1934.7142
$ws.Range("N122").Value = -7990.6

$ws.Range("H125").Value = 1000.1429
$ws.Range("I125").Value = 1010.75
$ws.Range("J125").Value = 986
$ws.Range("K125").Value = 9096.75
$ws.Range("L125").Value = 8874
$ws.Range("M125").Value = -6636.75
$ws.Range("N125").Value = -13794

$ws.Range("H132").Value = 2011.625
$ws.Range("I132").Value = 2015.6666
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 6046.9998
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -3516.9998
$ws.Range("N132").Value = -11058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1148.75
$ws.Range("I74").Value = 1148.75
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1148.75
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -274.75

$ws.Range("H77").Value = 1148.75
$ws.Range("I77").Value = 1148.75
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5743.75
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1375.75

$ws.Range("H132").Value = 1329.5
$ws.Range("I132").Value = 1333.1428
$ws.Range("J132").Value = 1304
$ws.Range("K132").Value = 3999.4284
$ws.Range("L132").Value = 3912
$ws.Range("M132").Value = -1469.4284
$ws.Range("N132").Value = -8972

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 9000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 9000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 9000
$ws.Range("N15").Value = -9454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 199
$ws.Range("I8").Value = 199
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 597
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -458

$ws.Range("H12").Value = 244.89285
$ws.Range("I12").Value = 169
$ws.Range("J12").Value = 275.25
$ws.Range("K12").Value = 507
$ws.Range("L12").Value = 825.75
$ws.Range("M12").Value = -334
$ws.Range("N12").Value = -1171.75

$ws.Range("H32").Value = 15008.5
$ws.Range("I32").Value = 2269
$ws.Range("J32").Value = 16828.428
$ws.Range("K32").Value = 6807
$ws.Range("L32").Value = 50485.284
$ws.Range("M32").Value = -6524
$ws.Range("N32").Value = -51051.284

$ws.Range("H95").Value = 7000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 7000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 21000
$ws.Range("N95").Value = -25118

$ws.Range("H107").Value = 1369.6666
$ws.Range("I107").Value = 445.5
$ws.Range("J107").Value = 1831.75
$ws.Range("K107").Value = 1336.5
$ws.Range("L107").Value = 5495.25
$ws.Range("M107").Value = 583.5
$ws.Range("N107").Value = -9335.25

$ws.Range("H139").Value = 1898.75
$ws.Range("I139").Value = 1898.75
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5696.25
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -556.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 349
$ws.Range("I2").Value = 6.142857
$ws.Range("J2").Value = 1549
$ws.Range("K2").Value = 6.142857
$ws.Range("L2").Value = 1549
$ws.Range("M2").Value = 106.857143
$ws.Range("N2").Value = -1775

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2471
$ws.Range("N22").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 37599.75
$ws.Range("I122").Value = 1340.2174
$ws.Range("J122").Value = 204393.6
$ws.Range("K122").Value = 4020.6522
$ws.Range("L122").Value = 613180.8
$ws.Range("M122").Value = -1570.6522
$ws.Range("N122").Value = -618080.8

$ws.Range("H126").Value = 6013.3335
$ws.Range("I126").Value = 6012
$ws.Range("J126").Value = 6014
$ws.Range("K126").Value = 18036
$ws.Range("L126").Value = 18042
$ws.Range("M126").Value = -15566
$ws.Range("N126").Value = -22982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3908.9092
$ws.Range("I22").Value = 1388.8889
$ws.Range("J22").Value = 15249
$ws.Range("K22").Value = 1388.8889
$ws.Range("L22").Value = 15249
$ws.Range("M22").Value = -1093.8889
$ws.Range("N22").Value = -15839

$ws.Range("H24").Value = 11000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 11000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 11000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -11686

$ws.Range("H27").Value = 3908.9092
$ws.Range("I27").Value = 1388.8889
$ws.Range("J27").Value = 15249
$ws.Range("K27").Value = 1388.8889
$ws.Range("L27").Value = 15249
$ws.Range("M27").Value = -1281.8889
$ws.Range("N27").Value = -15463

$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -812

$ws.Range("H61").Value = 4331.3076
$ws.Range("I61").Value = 4300.636
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 4300.636
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -4098.636
$ws.Range("N61").Value = -4904

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H93").Value = 1477.5555
$ws.Range("I93").Value = 1291.6154
$ws.Range("J93").Value = 1961
$ws.Range("K93").Value = 1291.6154
$ws.Range("L93").Value = 1961
$ws.Range("M93").Value = -43.61539999999991
$ws.Range("N93").Value = -4457

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H113").Value = 4331.3076
$ws.Range("I113").Value = 4300.636
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 4300.636
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -2130.636
$ws.Range("N113").Value = -8840

$ws.Range("H133").Value = 49999.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49999.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49999.332
$ws.Range("N133").Value = -55059.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33345
$ws.Range("I54").Value = 33345
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 33345
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -32825
$ws.Range("N54").ClearContents()

$ws.Range("H81").Value = 6078.364
$ws.Range("I81").Value = 6086.2
$ws.Range("J81").Value = 6000
$ws.Range("K81").Value = 12172.4
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -11111.4
$ws.Range("N81").Value = -14122

$ws.Range("H84").Value = 6078.364
$ws.Range("I84").Value = 6086.2
$ws.Range("J84").Value = 6000
$ws.Range("K84").Value = 60862
$ws.Range("L84").Value = 60000
$ws.Range("M84").Value = -55558
$ws.Range("N84").Value = -70608

$ws.Range("H122").Value = 1858.75
$ws.Range("I122").Value = 1711.1111
$ws.Range("J122").Value = 2301.6667
$ws.Range("K122").Value = 5133.3333
$ws.Range("L122").Value = 6905.000100000001
$ws.Range("M122").Value = -2683.3333
$ws.Range("N122").Value = -11805.0001

$ws.Range("H136").Value = 1860.3334
$ws.Range("I136").Value = 1040.5
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 3121.5
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -571.5
$ws.Range("N136").Value = -15600

